# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-15 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-12-16 Monday", 2)

# Update the division-problem table. The table's structure (20 rows x 5
# columns, with rows 2-4, 6-8, 10-12, 14-16, 18-20 blank) is unchanged;
# only the text of the populated cells changes. Target each cell by its
# (row, column) position so the edit is unambiguous even where the new
# text for one cell matches the old text that used to live in another
# cell (e.g. row 9).
$t = $d.Tables.Item(1)

# Row 1 (Word table row 1)
$t.Cell(1, 1).Range.Text = "162÷9=18, 0"
$t.Cell(1, 2).Range.Text = "143÷4=35, 3"
$t.Cell(1, 3).Range.Text = "224÷8=28, 0"
$t.Cell(1, 4).Range.Text = "484÷8=60, 4"
$t.Cell(1, 5).Range.Text = "336÷9=37, 3"

# Row 2 (Word table row 5)
$t.Cell(5, 1).Range.Text = "698÷5=139, 3"
$t.Cell(5, 2).Range.Text = "917÷3=305, 2"
$t.Cell(5, 3).Range.Text = "929÷8=116, 1"
$t.Cell(5, 4).Range.Text = "340÷6=56, 4"
$t.Cell(5, 5).Range.Text = "855÷3=285, 0"

# Row 3 (Word table row 9)
$t.Cell(9, 1).Range.Text = "457÷2=228, 1"
$t.Cell(9, 2).Range.Text = "445÷9=49, 4"
$t.Cell(9, 3).Range.Text = "509÷4=127, 1"
$t.Cell(9, 4).Range.Text = "578÷2=289, 0"
$t.Cell(9, 5).Range.Text = "977÷6=162, 5"

# Row 4 (Word table row 13)
$t.Cell(13, 1).Range.Text = "613÷5=122, 3"
$t.Cell(13, 2).Range.Text = "684÷8=85, 4"
$t.Cell(13, 3).Range.Text = "848÷3=282, 2"
$t.Cell(13, 4).Range.Text = "345÷7=49, 2"
$t.Cell(13, 5).Range.Text = "500÷8=62, 4"

# Row 5 (Word table row 17)
$t.Cell(17, 1).Range.Text = "976÷8=122, 0"
$t.Cell(17, 2).Range.Text = "891÷8=111, 3"
$t.Cell(17, 3).Range.Text = "472÷5=94, 2"
$t.Cell(17, 4).Range.Text = "245÷4=61, 1"
$t.Cell(17, 5).Range.Text = "339÷6=56, 3"
